# edit.ps1 - PowerPoint COM-interop script (iron_native run_com)
#
# Reproduces the authored change:
#   "added couts in count lambdas and hrvs"
#
# Concretely, on the single slide of the deck:
#   1. The "Title 4" shape (id=8, right-aligned author/affiliation text box
#      in the top-right corner) is made taller (its a:ext/cy grows) and two
#      new right-aligned paragraphs are appended to its text body:
#         "Research Group: Theoretical Physics"
#         <empty paragraph>
#      (its y-offset also shifts by 1 EMU as a natural side effect of the
#      resize in the source file).
#   2. The "TextBox 58" shape (id=59, the HKMM theorem / discrete spacetime
#      text box) is nudged to the right (its a:off/x increases).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "Title 4" (author / supervisor / research-group textbox) ----------
$title = $s.Shapes.Item(3)

# Resize/reposition to match the new taller text box.
$title.Top    = 5.690629921259842   # 72271 EMU  (was 72272 EMU)
$title.Height = 207.32055118110236  # 2632971 EMU (was 1441724 EMU)

# Append the two new paragraphs at the end of the existing text, preserving
# the existing runs/paragraphs untouched. The new paragraphs naturally
# inherit the right alignment (algn="r") from the preceding paragraph.
$titleText = $title.TextFrame.TextRange
$titleText.InsertAfter("`rResearch Group: Theoretical Physics`r")

# --- 2. "TextBox 58" (HKMM theorem / discreteness text box) ---------------
$theorem = $s.Shapes.Item(28)
$theorem.Left = 44.746230  # 568277 EMU (was 344871 EMU)
